$wb = $excel.ActiveWorkbook
$wsCodebook = $wb.Worksheets.Item("Codebook")
$wsDict = $wb.Worksheets.Item("Dictionary Mapping")

# New rows appended to the "Dictionary Mapping" sheet (rows 14-26)
$wsDict.Range("A14").Value = "??vic"
$wsDict.Range("F14").Value = "chebi:176783"
$wsDict.Range("H14").Value = "sio:SIO_000068"
$wsDict.Range("I14").Value = "??blood"

$wsDict.Range("A15").Value = "LBXVIC"
$wsDict.Range("B15").Value = "sio:SIO_001088"
$wsDict.Range("C15").Value = "??blood"
$wsDict.Range("D15").Value = "uo:0000165"
$wsDict.Range("I15").Value = "??vic"

$wsDict.Range("A16").Value = "??ucm"
$wsDict.Range("F16").Value = "chebi: 33007"
$wsDict.Range("I16").Value = "??urine"

$wsDict.Range("A17").Value = "URXUCM"
$wsDict.Range("B17").Value = "sio:SIO_"
$wsDict.Range("C17").Value = "??urine"
$wsDict.Range("I17").Value = "??ucm"

$wsDict.Range("A18").Value = "??ins"
$wsDict.Range("F18").Value = "chebi:145180"
$wsDict.Range("I18").Value = "??isn"

$wsDict.Range("A19").Value = "LBXIN"
$wsDict.Range("B19").Value = "sio:SIO_"
$wsDict.Range("C19").Value = "??blood"
$wsDict.Range("I19").Value = "??blood"

$wsDict.Range("A20").Value = "??uio"
$wsDict.Range("F20").Value = "chebi:33115"
$wsDict.Range("I20").Value = "??uio"

$wsDict.Range("A21").Value = "WTSA2YR"
$wsDict.Range("B21").Value = "sio:SIO_"
$wsDict.Range("C21").Value = "??urine"
$wsDict.Range("I21").Value = "??urine"

$wsDict.Range("A22").Value = "??vid"
$wsDict.Range("F22").Value = "chebi:27300"
$wsDict.Range("I22").Value = "??vid"

$wsDict.Range("A23").Value = "LBXVIDMS"
$wsDict.Range("B23").Value = "sio:SIO"
$wsDict.Range("C23").Value = "??blood"
$wsDict.Range("I23").Value = "??blood"

$wsDict.Range("A24").Value = "??hepb"
$wsDict.Range("I24").Value = "??hepb"

$wsDict.Range("A25").Value = "LBXHBC"
$wsDict.Range("B25").Value = "sio:SIO_"
$wsDict.Range("C25").Value = "??blood"

$wsDict.Range("A26").Value = "LBXBHS"
$wsDict.Range("C26").Value = "??blood"

# Update the active sheet / selection on each affected sheet to match the
# final view state: "Codebook" is no longer the selected tab, and
# "Dictionary Mapping" becomes the active tab with a new selection.
$wsCodebook.Activate()
$wsCodebook.Range("C11").Select()

$wsDict.Activate()
$wsDict.Range("D29").Select()
